$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Tabelle1")

$xlPasteFormats = -4122
$xlCenter = -4108

# --- New column E width (best achievable approximation of 45.5546875) ---
$ws.Columns("E").ColumnWidth = 44.7

# --- Row 18: extend header row with new "Tatsächliches Ergebnis" / "Bewertung" columns ---
# Copy the shaded header style from D18 onto the two new header cells, then set their text.
$ws.Range("D18").Copy()
$ws.Range("E18:F18").PasteSpecial($xlPasteFormats)
$ws.Range("E18").Value = "Tatsächliches Ergebnis"
$ws.Range("F18").Value = "Bewertung"

# --- Row 19 (Testcase 1): add Actual Result (copy of expected) + Bewertung "OK" ---
$ws.Rows(19).RowHeight = 22.2

$ws.Range("D19").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)
$ws.Range("E19").Value = $ws.Range("D19").Value2

$ws.Range("B19").Copy()
$ws.Range("F19").PasteSpecial($xlPasteFormats)
$ws.Range("F19").Value = "OK"

# --- Row 20 (Testcase 2) ---
$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial($xlPasteFormats)
$ws.Range("E20").Value = $ws.Range("D20").Value2

$ws.Range("B20").Copy()
$ws.Range("F20").PasteSpecial($xlPasteFormats)
$ws.Range("F20").Value = "OK"

# --- Row 21 (Testcase 3) - filter feature not implemented, Actual Result cell keeps default style ---
$ws.Range("E21").Value = "Nicht Implementiert"

$ws.Range("B21").Copy()
$ws.Range("F21").PasteSpecial($xlPasteFormats)
$ws.Range("F21").Value = "H"

# --- Row 22 (Testcase 4) ---
$ws.Range("D22").Copy()
$ws.Range("E22").PasteSpecial($xlPasteFormats)
$ws.Range("E22").Value = $ws.Range("D22").Value2

$ws.Range("B22").Copy()
$ws.Range("F22").PasteSpecial($xlPasteFormats)
$ws.Range("F22").Value = "OK"

# --- Row 24 (new): Fazit / conclusion ---
$ws.Rows(24).RowHeight = 43.2

$ws.Range("B18").Copy()
$ws.Range("B24").PasteSpecial($xlPasteFormats)
$ws.Range("B24").Value = "Fazit"

$ws.Range("C7").Copy()
$ws.Range("C24").PasteSpecial($xlPasteFormats)
$ws.Range("C24").Value = "Wegen Zeitmangel konnte die Filterfunktion nicht mehr vollständig implementiert werden. Alle anderen Funktionen sind vorhanden"

$ws.Application.CutCopyMode = $false

# --- Match final selection state ---
$ws.Range("C24").Select()
